$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Fill-ReflectionCell($row, $text) {
    $cell = $t.Cell($row, 2)
    $escaped = $text -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>' + $escaped + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $cell.Range.InsertXML($xml)

    $count = $cell.Range.Paragraphs.Count
    for ($i = $count; $i -ge 2; $i--) {
        $p = $cell.Range.Paragraphs.Item($i)
        $p.Range.Delete()
    }
}

Fill-ReflectionCell 2 'If I could redo this assignment, I would allocate more time to experimenting with advanced text representation methods such as word embeddings (Word2Vec, GloVe) and transformer-based embeddings (BERT) for both the topic modelling and sentiment classification tasks. For Task 1, I would explore BERTopic as an alternative to LDA to capture semantic relationships better. For Task 2, I would implement cross-validation earlier in the process rather than relying solely on a single train-test split, and I would explore more sophisticated handling of the class imbalance in the sentiment dataset using techniques like SMOTE or class weighting.'

Fill-ReflectionCell 3 "Through this assignment, I gained practical experience in end-to-end NLP pipelines, from data cleaning and text preprocessing to model building and evaluation. I learned how to use scikit-learn's LDA implementation for topic modelling, including hyperparameter tuning with grid search. I also developed skills in comparing multiple classification algorithms (Logistic Regression, Naive Bayes, Random Forest) across different text representations (TF-IDF vs BoW). Additionally, I improved my data visualization skills using matplotlib and seaborn for presenting model results effectively."

Fill-ReflectionCell 5 "I demonstrated innovation by designing a systematic approach to compare multiple models across different text representations, creating a comprehensive evaluation framework. For the topic modelling task, I creatively combined perplexity-based topic number selection with extensive hyperparameter tuning across alpha, beta, and learning decay parameters. For sentiment classification, I innovatively combined title and review text to create a richer feature set, and customized the stopword list to preserve sentiment-carrying words like 'not' and 'never' that are typically removed."

Fill-ReflectionCell 6 'I applied critical thinking throughout the assignment by justifying each data preprocessing decision with clear rationale. For example, I analysed why removing neutral reviews (3-star ratings) was necessary for binary sentiment classification, and evaluated the trade-offs of different text cleaning steps. I critically compared model performance using multiple evaluation metrics (accuracy, precision, recall, F1-score, ROC-AUC) rather than relying on a single metric, recognizing that each metric reveals different aspects of model quality. I also analysed feature importance to validate that the models learned meaningful sentiment patterns.'

Write-Output "done"
